# FIX: Agrega validacion en el caso de parametros nulos en el archivo de
# configuracion por filas vacias y actualiza expresion regular para detectar
# el nombre de los diplomas.

$wb = $excel.ActiveWorkbook

$rutas = $wb.Worksheets.Item("Rutas")
$documentos = $wb.Worksheets.Item("Documentos")

# --- Hoja "Rutas": se intercambian las rutas de entrada/salida ---
$rutas.Range("B2").Value = "C:\Diplomasporprocesar"
$rutas.Range("B3").Value = "C:\Septiembre"

# --- Hoja "Documentos": se agrega una fila nueva con la expresion regular
#     actualizada para detectar el nombre del diploma, conservando la
#     expresion anterior como respaldo en la fila siguiente ---
$documentos.Rows.Item(3).Insert()

$documentos.Range("A3").Value = "DIPLOMA"
$documentos.Range("B3").Value = "Diploma"
$documentos.Range("C3").Value = "DEL\s+A`u{00D1}O\s+\d{4}\.([^\n]+)"
$documentos.Range("D3").Value = "DIPLOMA"

$documentos.Range("C3").Style = $documentos.Range("C2").Style

# Activar la hoja Documentos, que es la que queda visible al guardar
$documentos.Activate()
